# Insert a new data row at row 96 (Vega Central Mapocho de Santiago - Arveja Verde).
# Existing rows 96-137 shift down to 97-138, preserving their data untouched.
# The freshly inserted row 96 gets a new market-report entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 96..137 down by one to make room for the new record.
$ws.Rows("96:96").Insert()

# Populate the newly inserted row 96 with the new observation.
$ws.Range("A96").Value2 = 9
$ws.Range("B96").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C96").Value2 = "Metropolitana"
$ws.Range("D96").Value2 = 44825
$ws.Range("E96").Value2 = 13
$ws.Range("F96").Value2 = 100112022
$ws.Range("G96").Value2 = "Arveja Verde"
$ws.Range("H96").Value2 = "Perfection"
$ws.Range("I96").Value2 = "Primera"
$ws.Range("J96").Value2 = 25
$ws.Range("K96").Value2 = 30000
$ws.Range("L96").Value2 = 30000
$ws.Range("M96").Value2 = 30000
$ws.Range("N96").Value2 = "$/malla 25 kilos"
$ws.Range("O96").Value2 = "Provincia de Limarí"
$ws.Range("P96").Value2 = 1200
$ws.Range("Q96").Value2 = 25
$ws.Range("R96").Value2 = "Hortaliza"
